$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New "Remark" header column (G1), formatted like the other header cells.
$ws.Range("G1").Value = "Remark"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# 2. Row 2 (Sugar 5kg) - Total Amount corrected.
$ws.Range("E2").Value = 20000

# 3. Row 3 becomes the "Chicken" line with updated quantities/pricing and a remark.
$ws.Range("A3").Value = "Chicken"
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 1000
$ws.Range("E3").Value = 30000
$ws.Range("F3").Value = "N/A"

$ws.Range("G3").Value = "we only can give 30kg"
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# 4. New row 4 - "Beef" line, formatted like row 3.
$ws.Range("A3:F3").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)

$ws.Range("A4").Value = "Beef"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 1500
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 30000
$ws.Range("F4").Value = "Available"

# 5. Column widths: Status column widened, new Remark column sized.
#    (ColumnWidth is quantized to the sheet's pixel grid on write, same as
#    real Excel, so these inputs are chosen to land on the nearest
#    achievable stored width to the target 20.43 / 36.86.)
$ws.Columns.Item(6).ColumnWidth = 19.6
$ws.Columns.Item(7).ColumnWidth = 35.95
